# Updated CHE model - 2025-09-06 17:10
# Applies the "fuels" sheet restructuring: coal/gas/oil -> COAL/GAS/OIL,
# split the "wind" fuel into "windon"/"windoff", backfill the Set (column B)
# for every NRG/DEM/ENV row, and reword the two CHE generation descriptions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fuels")

# Insert a new row above the old "ELC" row (row 15) so the "wind" row (14)
# can be split into two rows: windon (14) and windoff (new row 15).
$ws.Rows.Item(15).Insert()

# --- rename the three fossil fuels to their uppercase codes ---
$ws.Range("C6").Value = "COAL"
$ws.Range("C7").Value = "GAS"
$ws.Range("C12").Value = "OIL"

# --- split "wind" into onshore / offshore rows ---
$ws.Range("C14").Value = "windon"
$ws.Range("D14").Value = "wind onshore"

$ws.Range("C15").Value = "windoff"
$ws.Range("D15").Value = "wind offshore"
$ws.Range("E15").Value = "twh"

# --- backfill column B (Set) for rows that previously left it blank ---
$ws.Range("B6").Value = "NRG"
$ws.Range("B7").Value = "NRG"
$ws.Range("B8").Value = "NRG"
$ws.Range("B9").Value = "NRG"
$ws.Range("B10").Value = "NRG"
$ws.Range("B11").Value = "NRG"
$ws.Range("B12").Value = "NRG"
$ws.Range("B13").Value = "NRG"
$ws.Range("B14").Value = "NRG"
$ws.Range("B15").Value = "NRG"
$ws.Range("B16").Value = "NRG"
$ws.Range("B23").Value = "ENV"
$ws.Range("B24").Value = "ENV"

# --- reword the CHE generation descriptions ---
$ws.Range("O5").Value = "solar electricity generation"
$ws.Range("O6").Value = "onshore wind electricity generation"

# --- restore the active selection Excel recorded after the edit ---
$ws.Range("C13").Select()
